$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 17
$ws.Range("H17").Value = 59757.06
$ws.Range("J17").Value = 59757.06
$ws.Range("L17").Value = 179271.18
$ws.Range("N17").Value = -179607.18

# row 100
$ws.Range("H100").Value = 3676.5386
$ws.Range("J100").Value = 3998.3333
$ws.Range("L100").Value = 3998.3333
$ws.Range("N100").Value = -5080.3333

# row 113
$ws.Range("H113").Value = 2555.2307
$ws.Range("I113").Value = 2082.8572
$ws.Range("J113").Value = 3106.3333
$ws.Range("K113").Value = 2082.8572
$ws.Range("L113").Value = 3106.3333
$ws.Range("M113").Value = 1171.1428
$ws.Range("N113").Value = -9614.3333

# row 116
$ws.Range("H116").Value = 12502694
$ws.Range("I116").Value = 28573314
$ws.Range("J116").Value = 3322.111
$ws.Range("K116").Value = 28573314
$ws.Range("L116").Value = 3322.111
$ws.Range("M116").Value = -28569872
$ws.Range("N116").Value = -10206.111

$ws = $wb.Worksheets.Item("ARM")
# row 6
$ws.Range("H6").Value = 1000
$ws.Range("I6").Value = 1000
$ws.Range("K6").Value = 1000
$ws.Range("M6").Value = -827

# row 56
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# row 12
$ws.Range("H12").Value = 32633.334
$ws.Range("I12").Value = 17266.666
$ws.Range("J12").Value = 48000
$ws.Range("K12").Value = 17266.666
$ws.Range("L12").Value = 48000
$ws.Range("M12").Value = -17098.666
$ws.Range("N12").Value = -48336

# row 20
$ws.Range("H20").Value = 27318.281
$ws.Range("I20").Value = 38105.445
$ws.Range("J20").Value = 3047.1667
$ws.Range("K20").Value = 38105.445
$ws.Range("L20").Value = 3047.1667
$ws.Range("M20").Value = -37858.445
$ws.Range("N20").Value = -3541.1667

# row 107
$ws.Range("H107").Value = 18712.5
$ws.Range("I107").Value = 66062.5
$ws.Range("J107").Value = 2929.1667
$ws.Range("K107").Value = 66062.5
$ws.Range("L107").Value = 2929.1667
$ws.Range("M107").Value = -64142.5
$ws.Range("N107").Value = -6769.1667

$ws = $wb.Worksheets.Item("CRP")
# row 94
$ws.Range("H94").Value = 13889.5
$ws.Range("I94").Value = 1500
$ws.Range("J94").Value = 20084.25
$ws.Range("K94").Value = 1500
$ws.Range("L94").Value = 20084.25
$ws.Range("M94").Value = -1049
$ws.Range("N94").Value = -20986.25

# row 132
$ws.Range("H132").Value = 323260.88
$ws.Range("I132").Value = 436965.34
$ws.Range("J132").Value = 2821
$ws.Range("K132").Value = 1310896.02
$ws.Range("L132").Value = 8463
$ws.Range("M132").Value = -1308366.02
$ws.Range("N132").Value = -13523

$ws = $wb.Worksheets.Item("CUL")
# row 47
$ws.Range("H47").Value = 621.8
$ws.Range("I47").Value = 509.5
$ws.Range("J47").Value = 696.6667
$ws.Range("K47").Value = 1528.5
$ws.Range("L47").Value = 2090.0001
$ws.Range("M47").Value = -1097.5
$ws.Range("N47").Value = -2952.0001

# row 80
$ws.Range("H80").Value = 12990
$ws.Range("J80").Value = 12128.571
$ws.Range("L80").Value = 36385.713
$ws.Range("N80").Value = -38257.713

# row 83
$ws.Range("H83").Value = 12990
$ws.Range("J83").Value = 12128.571
$ws.Range("L83").Value = 109157.139
$ws.Range("N83").Value = -118517.139

# row 92
$ws.Range("H92").Value = 683.1667
$ws.Range("J92").Value = 683.1667
$ws.Range("L92").Value = 2049.5001
$ws.Range("N92").Value = -4545.5001

# row 97
$ws.Range("H97").Value = 999.6667
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 999.6667
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 2999.0001
$ws.Range("N97").Value = -3991.0001
$ws.Range("M97").ClearContents()

# row 98
$ws.Range("H98").Value = 1484.6364
$ws.Range("I98").Value = 2423.8
$ws.Range("J98").Value = 702
$ws.Range("K98").Value = 7271.400000000001
$ws.Range("L98").Value = 2106
$ws.Range("M98").Value = -5773.400000000001
$ws.Range("N98").Value = -5102

# row 107
$ws.Range("H107").Value = 355.16666
$ws.Range("I107").Value = 380.3913
$ws.Range("K107").Value = 1141.1739
$ws.Range("M107").Value = 778.8261

$ws = $wb.Worksheets.Item("GSM")
# row 102
$ws.Range("H102").Value = 2208.1035
$ws.Range("I102").Value = 1870.3182
$ws.Range("K102").Value = 1870.3182
$ws.Range("M102").Value = -248.3181999999999

# row 113
$ws.Range("H113").Value = 1493
$ws.Range("I113").Value = 1238
$ws.Range("J113").Value = 1668.3125
$ws.Range("K113").Value = 1238
$ws.Range("L113").Value = 1668.3125
$ws.Range("M113").Value = 932
$ws.Range("N113").Value = -6008.3125

$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 3674.8572
$ws.Range("I7").Value = 3648
$ws.Range("J7").Value = 3728.5715
$ws.Range("K7").Value = 3648
$ws.Range("L7").Value = 3728.5715
$ws.Range("M7").Value = -3536
$ws.Range("N7").Value = -3952.5715

# row 16
$ws.Range("H16").Value = 1046.9259
$ws.Range("I16").Value = 1074.6
$ws.Range("J16").Value = 701
$ws.Range("K16").Value = 1074.6
$ws.Range("L16").Value = 701
$ws.Range("M16").Value = -904.5999999999999
$ws.Range("N16").Value = -1041

# row 22
$ws.Range("H22").Value = 2000
$ws.Range("I22").Value = 2000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2000
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1705
$ws.Range("N22").ClearContents()

# row 27
$ws.Range("H27").Value = 2000
$ws.Range("I27").Value = 2000
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 2000
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -1893
$ws.Range("N27").ClearContents()

# row 40
$ws.Range("H40").Value = 3397.5557
$ws.Range("I40").Value = 3329.6667
$ws.Range("J40").Value = 3533.3333
$ws.Range("K40").Value = 3329.6667
$ws.Range("L40").Value = 3533.3333
$ws.Range("M40").Value = -3193.6667
$ws.Range("N40").Value = -3805.3333

# row 122
$ws.Range("H122").Value = 13339173
$ws.Range("I122").Value = 4885.7144
$ws.Range("J122").Value = 25006674
$ws.Range("K122").Value = 14657.1432
$ws.Range("L122").Value = 75020022
$ws.Range("M122").Value = -12207.1432
$ws.Range("N122").Value = -75024922

# row 126
$ws.Range("H126").Value = 3674.8572
$ws.Range("I126").Value = 3648
$ws.Range("J126").Value = 3728.5715
$ws.Range("K126").Value = 10944
$ws.Range("L126").Value = 11185.7145
$ws.Range("M126").Value = -8474
$ws.Range("N126").Value = -16125.7145

$ws = $wb.Worksheets.Item("WVR")
# row 11
$ws.Range("H11").Value = 20000
$ws.Range("J11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("N11").Value = -10284

# row 96
$ws.Range("H96").Value = 1850
$ws.Range("I96").Value = 1966.6666
$ws.Range("J96").Value = 1500
$ws.Range("K96").Value = 1966.6666
$ws.Range("L96").Value = 1500
$ws.Range("M96").Value = -593.6666
$ws.Range("N96").Value = -4246

# row 126
$ws.Range("H126").Value = 5254
$ws.Range("I126").Value = 7492.6665
$ws.Range("J126").Value = 776.6667
$ws.Range("K126").Value = 22477.9995
$ws.Range("L126").Value = 2330.0001
$ws.Range("M126").Value = -20007.9995
$ws.Range("N126").Value = -7270.0001

# row 132
$ws.Range("H132").Value = 1743.75
$ws.Range("I132").Value = 1018.2432
$ws.Range("J132").Value = 3533.3333
$ws.Range("K132").Value = 3054.7296
$ws.Range("L132").Value = 10599.9999
$ws.Range("M132").Value = -524.7296000000001
$ws.Range("N132").Value = -15659.9999

# row 136
$ws.Range("H136").Value = 1650.931
$ws.Range("I136").Value = 1540.7084
$ws.Range("J136").Value = 2180
$ws.Range("K136").Value = 4622.1252
$ws.Range("L136").Value = 6540
$ws.Range("M136").Value = -2072.1252
$ws.Range("N136").Value = -11640
